$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("A1").Value = "Url"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "dbExcel"
$ws.Range("D1").Value = "WebExcel"

# --- Row 2 data ---
$ws.Range("B2").Value = @'
MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"] and d.er_status In ["Positive"]and d.pr_status In ["Positive"] 
return DISTINCT ss.study_subject_id as `Case ID`,
   p.program_acronym as `Program Code`,
    p.program_id as Program_ID,
   s.study_acronym as `Arm`,
   ss.disease_subtype as `Diagnosis`,
   sf.grouped_recurrence_score AS `Recurrence Score`,
   d.tumor_size_group AS `tumor_size`,
   d.er_status AS `ER Status`,
   d.pr_status AS `PR Status`,
   demo.age_at_index AS `Age (years)`,
	demo.survival_time AS `Survival (days)`
'@

$ws.Range("C2").Value = "TC03_Bento_E2E_Select-Single-CaseDetail_Manifest.xlsx"
$ws.Range("D2").Value = "TC03_Bento_E2E_Select-Single-CaseDetail_WebData.xlsx"

# A2 becomes a hyperlink to the Bento QA site; reset to Normal first so the
# new Hyperlink style doesn't inherit the old wrap-text alignment.
$ws.Range("A2").Style = "Normal"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://bento-qa.bento-tools.org/", "", "", "https://bento-qa.bento-tools.org/") | Out-Null

# B2 keeps the wrapped-text style used by the long Cypher query.
$ws.Range("B2").WrapText = $true

# Row 2 shrinks to fit the shorter query text.
$ws.Rows.Item(2).RowHeight = 375

# Restore default view: selection moves to D7 (scrolled back to top-left A1).
$ws.Range("D7").Select() | Out-Null
